$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = '{"UserInfo":{"ID":"atpermission@iicon004.com","DisplayName":"Auto+Test+Permissions","TextAddress":"","Mobile":"","BusinessPhone":"","Active":true,"UserType":"1","IsUser":true,"AllowLoginIntoPC":true,"AllowLoginIntoInspectMobile":false,"AllowLoginIntoFleetMobile":false,"AllowMobileBarcodeScanning":false,"ContactType":"100","ManagerIID":"","Notes":"","AssignedWorkOrders":false,"EmailOptOut":false,"InspectEmailList":false,"TeamIntelligenceUser":false,"FOB":"","HourlyRate":-1,"LandingPage":"MapView.aspx","PreferredLanguage":"","IsAdvisor":false,"LoginVerifyType":"0","LocationIds":[],"DepartmentIds":[],"IID":"5AB99F30-4CC9-4C74-98C0-57A426A72114"},"Subscribe":{"$type":"FI.FIC.Contracts.DataObjects.BaseObject.SubscribeMessageByEmail,+FICIntfAdv","NeedSendMessages":[],"UserEmail":"","UserTextMessage":""},"Features":[{"Key":100,"Value":["0"]},{"Key":110,"Value":["0"]},{"Key":120,"Value":["0"]},{"Key":130,"Value":["0"]},{"Key":140,"Value":["0"]},{"Key":200,"Value":["0"]},{"Key":210,"Value":["0"]},{"Key":220,"Value":["0"]},{"Key":230,"Value":["0"]},{"Key":235,"Value":["0"]},{"Key":237,"Value":["0"]},{"Key":245,"Value":["0"]},{"Key":300,"Value":["0"]},{"Key":600,"Value":["99999"]},{"Key":601,"Value":["0"]},{"Key":602,"Value":["0"]},{"Key":800,"Value":["99999"]},{"Key":900,"Value":["0"]},{"Key":1000,"Value":["99999"]},{"Key":1100,"Value":["0"]},{"Key":1101,"Value":["0"]}],"Schedule":{"$type":"FI.FIC.EmailSchedule,+FICBLC","ScheduleItems":{"$type":"FI.FIC.EmailScheduleItem[],+FICBLC","$values":[]}}}'
$ws.Range("B3").Value = '{"UserInfo":{"ID":"atpermission@iicon004.com","DisplayName":"Auto+Test+Permissions","TextAddress":"","Mobile":"","BusinessPhone":"","Active":true,"UserType":"1","IsUser":true,"AllowLoginIntoPC":true,"AllowLoginIntoInspectMobile":false,"AllowLoginIntoFleetMobile":false,"AllowMobileBarcodeScanning":false,"ContactType":"100","ManagerIID":"","Notes":"","AssignedWorkOrders":false,"EmailOptOut":false,"InspectEmailList":false,"TeamIntelligenceUser":false,"FOB":"","HourlyRate":-1,"LandingPage":"MapView.aspx","PreferredLanguage":"","IsAdvisor":false,"LoginVerifyType":"0","LocationIds":[],"DepartmentIds":[],"IID":"5AB99F30-4CC9-4C74-98C0-57A426A72114"},"Subscribe":{"$type":"FI.FIC.Contracts.DataObjects.BaseObject.SubscribeMessageByEmail,+FICIntfAdv","NeedSendMessages":[],"UserEmail":"","UserTextMessage":""},"Features":[{"Key":100,"Value":["0"]},{"Key":110,"Value":["0"]},{"Key":120,"Value":["0"]},{"Key":130,"Value":["0"]},{"Key":140,"Value":["0"]},{"Key":200,"Value":["0"]},{"Key":210,"Value":["0"]},{"Key":220,"Value":["0"]},{"Key":230,"Value":["0"]},{"Key":235,"Value":["0"]},{"Key":237,"Value":["0"]},{"Key":245,"Value":["0"]},{"Key":300,"Value":["99999"]},{"Key":600,"Value":["0"]},{"Key":601,"Value":["0"]},{"Key":602,"Value":["0"]},{"Key":800,"Value":["99999"]},{"Key":900,"Value":["99999"]},{"Key":1000,"Value":["99999"]},{"Key":1100,"Value":["0"]},{"Key":1101,"Value":["0"]}],"Schedule":{"$type":"FI.FIC.EmailSchedule,+FICBLC","ScheduleItems":{"$type":"FI.FIC.EmailScheduleItem[],+FICBLC","$values":[]}}}'
$ws.Range("B4").Value = '{"UserInfo":{"ID":"atpermission@iicon004.com","DisplayName":"Auto+Test+Permissions","TextAddress":"","Mobile":"","BusinessPhone":"","Active":true,"UserType":"1","IsUser":true,"AllowLoginIntoPC":true,"AllowLoginIntoInspectMobile":false,"AllowLoginIntoFleetMobile":false,"AllowMobileBarcodeScanning":false,"ContactType":"100","ManagerIID":"","Notes":"","AssignedWorkOrders":false,"EmailOptOut":false,"InspectEmailList":false,"TeamIntelligenceUser":false,"FOB":"","HourlyRate":-1,"LandingPage":"MapView.aspx","PreferredLanguage":"","IsAdvisor":false,"LoginVerifyType":"0","LocationIds":[],"DepartmentIds":[],"IID":"5AB99F30-4CC9-4C74-98C0-57A426A72114"},"Subscribe":{"$type":"FI.FIC.Contracts.DataObjects.BaseObject.SubscribeMessageByEmail,+FICIntfAdv","NeedSendMessages":[],"UserEmail":"","UserTextMessage":""},"Features":[{"Key":100,"Value":["99999"]},{"Key":110,"Value":["0"]},{"Key":120,"Value":["0"]},{"Key":130,"Value":["0"]},{"Key":140,"Value":["0"]},{"Key":200,"Value":["99999"]},{"Key":210,"Value":["0"]},{"Key":220,"Value":["0"]},{"Key":230,"Value":["0"]},{"Key":235,"Value":["0"]},{"Key":237,"Value":["0"]},{"Key":245,"Value":["0"]},{"Key":300,"Value":["0"]},{"Key":600,"Value":["99999"]},{"Key":601,"Value":["0"]},{"Key":602,"Value":["0"]},{"Key":800,"Value":["0"]},{"Key":900,"Value":["0"]},{"Key":1000,"Value":["99999"]},{"Key":1100,"Value":["99999"]},{"Key":1101,"Value":["0"]}],"Schedule":{"$type":"FI.FIC.EmailSchedule,+FICBLC","ScheduleItems":{"$type":"FI.FIC.EmailScheduleItem[],+FICBLC","$values":[]}}}'

$ws.Range("C14").Select()
